# Website Changes tracker - add two new rows (Account / Course Contents)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 - Account
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "Account"
$ws.Range("C18").Value = "As of now hide Account button"

# Row 19 - Course Contents
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Course Contents"
$ws.Range("C19").Value = "Place the suitable image in the background with respect to courses"
$ws.Range("E19").Value = "For Exmaple, If we choose AWS, then background image should be AWS "

# Formatting to match the rest of the sheet: S.No/Activity/Comments columns
# are plain, Description (and the filler cell after it) wrap their text.
$ws.Range("A18:E19").WrapText = $false
$ws.Range("C18:D19").WrapText = $true
$ws.Range("E19").WrapText = $true

# Row 19 holds a long comment, so it needs to be considerably taller.
$ws.Rows.Item(19).RowHeight = 105

# Leave the cursor on the newly added block, matching where the author
# finished editing.
$null = $ws.Range("A18:E19").Select()
